# Apply "manual case resampling" re-run of the bootstrap results.
# Sheet "Full results": refreshed numeric estimates (same #NUM!/formula-free layout).
# Sheet "For plotting": header columns A/B swapped (Index, Outcome) and the
# whole 15-row body re-derived/re-ordered (grouped by Index then Outcome) with new values.

$wb = $excel.ActiveWorkbook
$wsFull = $wb.Worksheets.Item("Full results")
$wsPlot = $wb.Worksheets.Item("For plotting")

# --- "Full results": refreshed numeric estimates ---
$wsFull.Range("H2").Value = 0.588577598665504
$wsFull.Range("I2").Value = 0.278240951962408
$wsFull.Range("O2").Value = 0.411207907409411
$wsFull.Range("F3").Value = 0.60400366182116
$wsFull.Range("G3").Value = 0.310140672331196
$wsFull.Range("C4").Value = 0.638066291221125
$wsFull.Range("D4").Value = 0.361569411877401
$wsFull.Range("E4").Value = 0.999635703098526
$wsFull.Range("J4").Value = 0.36170117980571
$wsFull.Range("K4").Value = 0.310253697948053
$wsFull.Range("L4").Value = 0.0154316848992071
$wsFull.Range("M4").Value = 0.0495067276037014
$wsFull.Range("N4").Value = 0.32568538284726
$wsFull.Range("H5").Value = 0.871937806310243
$wsFull.Range("I5").Value = 0.0879278595061394
$wsFull.Range("O5").Value = 0.128085046372583
$wsFull.Range("F6").Value = 0.882606928066615
$wsFull.Range("G6").Value = 0.0933282465590723
$wsFull.Range("C7").Value = 0.89233621196557
$wsFull.Range("D7").Value = 0.107689997855455
$wsFull.Range("E7").Value = 1.00002620982103
$wsFull.Range("J7").Value = 0.107687175269545
$wsFull.Range("K7").Value = 0.0933258003795308
$wsFull.Range("L7").Value = 0.0106688421869644
$wsFull.Range("M7").Value = 0.0203978711030383
$wsFull.Range("N7").Value = 0.103994642566495
$wsFull.Range("H8").Value = 0.81427089119318
$wsFull.Range("I8").Value = 0.10902784328245
$wsFull.Range("O8").Value = 0.185982621113335
$wsFull.Range("F9").Value = 0.855337800394946
$wsFull.Range("G9").Value = 0.125183563805627
$wsFull.Range("C10").Value = 0.859701573261984
$wsFull.Range("D10").Value = 0.140609863257113
$wsFull.Range("E10").Value = 1.0003114365191
$wsFull.Range("J10").Value = 0.14056608463557
$wsFull.Range("K10").Value = 0.125144588036474
$wsFull.Range("L10").Value = 0.0410541222646858
$wsFull.Range("M10").Value = 0.0454165364777647
$wsFull.Range("N10").Value = 0.16619871030116
$wsFull.Range("H11").Value = 0.808662161004036
$wsFull.Range("I11").Value = 0.137461374945292
$wsFull.Range("O11").Value = 0.191621736009115
$wsFull.Range("F12").Value = 0.810197271584953
$wsFull.Range("G12").Value = 0.151187849552119
$wsFull.Range("C13").Value = 0.827898351653003
$wsFull.Range("D13").Value = 0.172452844834897
$wsFull.Range("E13").Value = 1.0003511964879
$wsFull.Range("J13").Value = 0.172392297681435
$wsFull.Range("K13").Value = 0.151134768609519
$wsFull.Range("L13").Value = 0.00153457203766954
$wsFull.Range("M13").Value = 0.0192294383276799
$wsFull.Range("N13").Value = 0.152669340647188
$wsFull.Range("H14").Value = 0.798854634973577
$wsFull.Range("I14").Value = 0.161234111127771
$wsFull.Range("O14").Value = 0.202339734090694
$wsFull.Range("F15").Value = 0.817704725280448
$wsFull.Range("G15").Value = 0.165582003926798
$wsFull.Range("C16").Value = 0.828573269432525
$wsFull.Range("D16").Value = 0.172924072006996
$wsFull.Range("E16").Value = 1.00149734143952
$wsFull.Range("J16").Value = 0.17266553142669
$wsFull.Range("K16").Value = 0.16533444054959
$wsFull.Range("L16").Value = 0.0188219080317273
$wsFull.Range("M16").Value = 0.0296742026640041
$wsFull.Range("N16").Value = 0.184156348581317

# --- "For plotting": header swap (A=Index, B=Outcome) ---
$wsPlot.Range("A1").Value = "Index"
$wsPlot.Range("B1").Value = "Outcome"

# --- "For plotting": rebuilt/reordered body rows 2-16 ---
$wsPlot.Range("A2").Value = "Sibcorr"
$wsPlot.Range("B2").Value = "education"
$wsPlot.Range("C2").Value = 0.36170117980571
$wsPlot.Range("D2").Value = 0.321846345256781
$wsPlot.Range("E2").Value = 0.401556014354639
$wsPlot.Range("A3").Value = "IOLIB"
$wsPlot.Range("B3").Value = "education"
$wsPlot.Range("C3").Value = 0.32568538284726
$wsPlot.Range("D3").Value = 0.284428294588826
$wsPlot.Range("E3").Value = 0.366942471105695
$wsPlot.Range("A4").Value = "IORAD"
$wsPlot.Range("B4").Value = "education"
$wsPlot.Range("C4").Value = 0.411207907409411
$wsPlot.Range("D4").Value = 0.373866616526509
$wsPlot.Range("E4").Value = 0.448549198292314
$wsPlot.Range("A5").Value = "Sibcorr"
$wsPlot.Range("B5").Value = "occupation"
$wsPlot.Range("C5").Value = 0.172392297681435
$wsPlot.Range("D5").Value = 0.140769610086204
$wsPlot.Range("E5").Value = 0.204014985276666
$wsPlot.Range("A6").Value = "IOLIB"
$wsPlot.Range("B6").Value = "occupation"
$wsPlot.Range("C6").Value = 0.152669340647188
$wsPlot.Range("D6").Value = 0.119133473603247
$wsPlot.Range("E6").Value = 0.18620520769113
$wsPlot.Range("A7").Value = "IORAD"
$wsPlot.Range("B7").Value = "occupation"
$wsPlot.Range("C7").Value = 0.191621736009115
$wsPlot.Range("D7").Value = 0.159968711699012
$wsPlot.Range("E7").Value = 0.223274760319218
$wsPlot.Range("A8").Value = "Sibcorr"
$wsPlot.Range("B8").Value = "income"
$wsPlot.Range("C8").Value = 0.14056608463557
$wsPlot.Range("D8").Value = 0.0661913249147279
$wsPlot.Range("E8").Value = 0.214940844356412
$wsPlot.Range("A9").Value = "IOLIB"
$wsPlot.Range("B9").Value = "income"
$wsPlot.Range("C9").Value = 0.16619871030116
$wsPlot.Range("D9").Value = 0.0979678299955207
$wsPlot.Range("E9").Value = 0.2344295906068
$wsPlot.Range("A10").Value = "IORAD"
$wsPlot.Range("B10").Value = "income"
$wsPlot.Range("C10").Value = 0.185982621113335
$wsPlot.Range("D10").Value = 0.115247247965152
$wsPlot.Range("E10").Value = 0.256717994261517
$wsPlot.Range("A11").Value = "Sibcorr"
$wsPlot.Range("B11").Value = "wealth"
$wsPlot.Range("C11").Value = 0.17266553142669
$wsPlot.Range("D11").Value = 0.0918797609897858
$wsPlot.Range("E11").Value = 0.253451301863595
$wsPlot.Range("A12").Value = "IOLIB"
$wsPlot.Range("B12").Value = "wealth"
$wsPlot.Range("C12").Value = 0.184156348581317
$wsPlot.Range("D12").Value = 0.106753355168802
$wsPlot.Range("E12").Value = 0.261559341993832
$wsPlot.Range("A13").Value = "IORAD"
$wsPlot.Range("B13").Value = "wealth"
$wsPlot.Range("C13").Value = 0.202339734090694
$wsPlot.Range("D13").Value = 0.131528250415873
$wsPlot.Range("E13").Value = 0.273151217765516
$wsPlot.Range("A14").Value = "Sibcorr"
$wsPlot.Range("B14").Value = "health_pc"
$wsPlot.Range("C14").Value = 0.107687175269545
$wsPlot.Range("D14").Value = 0.0644074612943271
$wsPlot.Range("E14").Value = 0.150966889244762
$wsPlot.Range("A15").Value = "IOLIB"
$wsPlot.Range("B15").Value = "health_pc"
$wsPlot.Range("C15").Value = 0.103994642566495
$wsPlot.Range("D15").Value = 0.0633626210390709
$wsPlot.Range("E15").Value = 0.14462666409392
$wsPlot.Range("A16").Value = "IORAD"
$wsPlot.Range("B16").Value = "health_pc"
$wsPlot.Range("C16").Value = 0.128085046372583
$wsPlot.Range("D16").Value = 0.0887273901220133
$wsPlot.Range("E16").Value = 0.167442702623153
